$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1238.1666
$ws.Range("I2").Value = 1306
$ws.Range("J2").Value = 899
$ws.Range("K2").Value = 1306
$ws.Range("L2").Value = 899
$ws.Range("M2").Value = -1193
$ws.Range("N2").Value = -1125
$ws.Range("H6").Value = 978.7692
$ws.Range("I6").Value = 652
$ws.Range("K6").Value = 1956
$ws.Range("M6").Value = -1844
$ws.Range("H32").Value = 18103.5
$ws.Range("I32").Value = 13766.6
$ws.Range("K32").Value = 13766.6
$ws.Range("M32").Value = -13440.6
$ws.Range("H100").Value = 3174.125
$ws.Range("I100").Value = 2444.8333
$ws.Range("K100").Value = 2444.8333
$ws.Range("M100").Value = -1903.8333
$ws.Range("H112").Value = 5817.778
$ws.Range("J112").Value = 6428.4585
$ws.Range("L112").Value = 19285.3755
$ws.Range("N112").Value = -21501.3755
$ws.Range("H118").Value = 989.61536
$ws.Range("I118").Value = 596.8182
$ws.Range("K118").Value = 1790.4546
$ws.Range("M118").Value = -133.4546
$ws.Range("H137").Value = 14712.826
$ws.Range("I137").Value = 1582
$ws.Range("K137").Value = 4746
$ws.Range("M137").Value = -2196
$ws.Range("H138").Value = 4332044.5
$ws.Range("I138").Value = 2816.3333
$ws.Range("J138").Value = 6496659
$ws.Range("K138").Value = 8448.999899999999
$ws.Range("L138").Value = 19489977
$ws.Range("M138").Value = -3308.999899999999
$ws.Range("N138").Value = -19500257

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20409.422
$ws.Range("I32").Value = 22508.32
$ws.Range("J32").Value = 5417.2856
$ws.Range("K32").Value = 22508.32
$ws.Range("L32").Value = 5417.2856
$ws.Range("M32").Value = -22221.32
$ws.Range("N32").Value = -5991.2856
$ws.Range("H61").Value = 5206.5713
$ws.Range("I61").Value = 3115.6
$ws.Range("J61").Value = 14499.777
$ws.Range("K61").Value = 3115.6
$ws.Range("L61").Value = 14499.777
$ws.Range("M61").Value = -2903.6
$ws.Range("N61").Value = -14923.777
$ws.Range("H74").Value = 5108.9355
$ws.Range("I74").Value = 1432.7916
$ws.Range("J74").Value = 17712.857
$ws.Range("K74").Value = 1432.7916
$ws.Range("L74").Value = 17712.857
$ws.Range("M74").Value = -558.7916
$ws.Range("N74").Value = -19460.857
$ws.Range("H77").Value = 5108.9355
$ws.Range("I77").Value = 1432.7916
$ws.Range("J77").Value = 17712.857
$ws.Range("K77").Value = 7163.958000000001
$ws.Range("L77").Value = 88564.285
$ws.Range("M77").Value = -2795.958000000001
$ws.Range("N77").Value = -97300.285
$ws.Range("H92").Value = 90000
$ws.Range("I92").Value = 90000
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 90000
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -87504
$ws.Range("N92").ClearContents()
$ws.Range("H132").Value = 2467.8853
$ws.Range("I132").Value = 2308.22
$ws.Range("J132").Value = 3193.6365
$ws.Range("K132").Value = 6924.66
$ws.Range("L132").Value = 9580.9095
$ws.Range("M132").Value = -4394.66
$ws.Range("N132").Value = -14640.9095
$ws.Range("H136").Value = 5206.5713
$ws.Range("I136").Value = 3115.6
$ws.Range("J136").Value = 14499.777
$ws.Range("K136").Value = 9346.799999999999
$ws.Range("L136").Value = 43499.331
$ws.Range("M136").Value = -6796.799999999999
$ws.Range("N136").Value = -48599.331

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3630.4285
$ws.Range("I99").Value = 3630.4285
$ws.Range("K99").Value = 3630.4285
$ws.Range("M99").Value = -2132.4285
$ws.Range("H107").Value = 1502.75
$ws.Range("I107").Value = 1502.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1502.75
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 417.25
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 5167.15
$ws.Range("I134").Value = 5760.2666
$ws.Range("J134").Value = 3387.8
$ws.Range("K134").Value = 17280.7998
$ws.Range("L134").Value = 10163.4
$ws.Range("M134").Value = -14745.7998
$ws.Range("N134").Value = -15233.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2502448
$ws.Range("I31").Value = 4168554.2
$ws.Range("K31").Value = 4168554.2
$ws.Range("M31").Value = -4168259.2
$ws.Range("H34").Value = 2502448
$ws.Range("I34").Value = 4168554.2
$ws.Range("K34").Value = 4168554.2
$ws.Range("M34").Value = -4168352.2
$ws.Range("H94").Value = 1121.7222
$ws.Range("I94").Value = 1488.7142
$ws.Range("J94").Value = 888.1818
$ws.Range("K94").Value = 1488.7142
$ws.Range("L94").Value = 888.1818
$ws.Range("M94").Value = -1037.7142
$ws.Range("N94").Value = -1790.1818
$ws.Range("H99").Value = 3492.3125
$ws.Range("I99").Value = 3198.75
$ws.Range("J99").Value = 3785.875
$ws.Range("K99").Value = 3198.75
$ws.Range("L99").Value = 3785.875
$ws.Range("M99").Value = -1700.75
$ws.Range("N99").Value = -6781.875
$ws.Range("H126").Value = 3492.3125
$ws.Range("I126").Value = 3198.75
$ws.Range("J126").Value = 3785.875
$ws.Range("K126").Value = 9596.25
$ws.Range("L126").Value = 11357.625
$ws.Range("M126").Value = -7126.25
$ws.Range("N126").Value = -16297.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3140.3704
$ws.Range("I68").Value = 1546.3334
$ws.Range("J68").Value = 3339.625
$ws.Range("K68").Value = 4639.0002
$ws.Range("L68").Value = 10018.875
$ws.Range("M68").Value = -3828.0002
$ws.Range("N68").Value = -11640.875
$ws.Range("H71").Value = 3140.3704
$ws.Range("I71").Value = 1546.3334
$ws.Range("J71").Value = 3339.625
$ws.Range("K71").Value = 13917.0006
$ws.Range("L71").Value = 30056.625
$ws.Range("M71").Value = -9861.000599999999
$ws.Range("N71").Value = -38168.625
$ws.Range("H131").Value = 2504.12
$ws.Range("J131").Value = 1952.6364
$ws.Range("L131").Value = 5857.9092
$ws.Range("N131").Value = -15937.9092

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H97").Value = 398
$ws.Range("I97").Value = 408
$ws.Range("J97").Value = 238
$ws.Range("K97").Value = 408
$ws.Range("L97").Value = 238
$ws.Range("M97").Value = 88
$ws.Range("N97").Value = -1230
$ws.Range("H122").Value = 14286804
$ws.Range("I122").Value = 725.8
$ws.Range("J122").Value = 50002000
$ws.Range("K122").Value = 2177.4
$ws.Range("L122").Value = 150006000
$ws.Range("M122").Value = 272.6000000000004
$ws.Range("N122").Value = -150010900
$ws.Range("H132").Value = 5086.4473
$ws.Range("I132").Value = 5423.926
$ws.Range("J132").Value = 4258.091
$ws.Range("K132").Value = 16271.778
$ws.Range("L132").Value = 12774.273
$ws.Range("M132").Value = -13741.778
$ws.Range("N132").Value = -17834.273
$ws.Range("H135").Value = 114056.29
$ws.Range("J135").Value = 114056.29
$ws.Range("L135").Value = 114056.29
$ws.Range("N135").Value = -124196.29
$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360
$ws.Range("H141").Value = 72500
$ws.Range("J141").Value = 72500
$ws.Range("L141").Value = 72500
$ws.Range("N141").Value = -82860

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1701.5834
$ws.Range("I22").Value = 1452.375
$ws.Range("K22").Value = 1452.375
$ws.Range("M22").Value = -1157.375
$ws.Range("H27").Value = 1701.5834
$ws.Range("I27").Value = 1452.375
$ws.Range("K27").Value = 1452.375
$ws.Range("M27").Value = -1345.375
$ws.Range("H40").Value = 2860.6785
$ws.Range("I40").Value = 2734.577
$ws.Range("K40").Value = 2734.577
$ws.Range("M40").Value = -2598.577
$ws.Range("H46").Value = 4178.2856
$ws.Range("I46").Value = 983
$ws.Range("J46").Value = 4710.8335
$ws.Range("K46").Value = 983
$ws.Range("L46").Value = 4710.8335
$ws.Range("M46").Value = -795
$ws.Range("N46").Value = -5086.8335
$ws.Range("H55").Value = 1826.0667
$ws.Range("I55").Value = 487.4
$ws.Range("J55").Value = 2495.4
$ws.Range("K55").Value = 487.4
$ws.Range("L55").Value = 2495.4
$ws.Range("M55").Value = -314.4
$ws.Range("N55").Value = -2841.4
$ws.Range("H132").Value = 3947.0527
$ws.Range("I132").Value = 3412.4614
$ws.Range("K132").Value = 10237.3842
$ws.Range("M132").Value = -7707.3842

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1050.6875
$ws.Range("I100").Value = 664.1111
$ws.Range("K100").Value = 1328.2222
$ws.Range("M100").Value = -787.2221999999999
$ws.Range("H107").Value = 788.1111
$ws.Range("I107").Value = 884
$ws.Range("K107").Value = 2652
$ws.Range("M107").Value = -732
$ws.Range("H126").Value = 2559.077
$ws.Range("I126").Value = 1933.8
$ws.Range("J126").Value = 4643.3335
$ws.Range("K126").Value = 5801.4
$ws.Range("L126").Value = 13930.0005
$ws.Range("M126").Value = -3331.4
$ws.Range("N126").Value = -18870.0005
$ws.Range("H132").Value = 12504232
$ws.Range("I132").Value = 12504232
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 37512696
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -37510166
$ws.Range("N132").ClearContents()
